$wb = $excel.ActiveWorkbook

# --- Overview sheet: roll up the new "Ready for handoff" status/date for
#     e1474d6f-3107-4f3e-b5ff-caeab908aca6.md (row 9) and refresh the
#     already-ready fbdf8557-9368-407b-a255-6254c559e860.md (row 10) to the
#     same latest handoff timestamp. ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B9").Value = "Ready for handoff"
$overview.Range("C9").Value = "Ready for handoff"
$overview.Range("D9").Value = "2016-19-17 18:19:19"
$overview.Range("D10").Value = "2016-19-17 18:19:19"

# --- zh-cn detail sheet: row 9 is e1474d6f, row 10 is fbdf8557 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C9").Value = "Ready for handoff"
$zhcn.Range("E9").Value = "2016-03-17 18:19:15"
$zhcn.Range("E10").Value = "2016-03-17 18:19:15"

# --- de-de detail sheet: row 9 is e1474d6f, row 10 is fbdf8557 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C9").Value = "Ready for handoff"
$dede.Range("E9").Value = "2016-03-17 18:19:19"
$dede.Range("E10").Value = "2016-03-17 18:19:19"
